$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 35
$ws.Range("F3").Value = 1292
$ws.Range("F4").Value = 12995
$ws.Range("F5").Value = 741
$ws.Range("F10").Value = 1887
$ws.Range("F11").Value = 41
$ws.Range("F13").Value = 4832
$ws.Range("F14").Value = 528
$ws.Range("F15").Value = 213
$ws.Range("F17").Value = 354
$ws.Range("F19").Value = 305
$ws.Range("F20").Value = 139
$ws.Range("F21").Value = 133
$ws.Range("F23").Value = 226
$ws.Range("F24").Value = 266
$ws.Range("F25").Value = 1315
$ws.Range("F26").Value = 348
# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 287
$ws.Range("F6").Value = 167
$ws.Range("F7").Value = 19
$ws.Range("F8").Value = 19
$ws.Range("F11").Value = 365
$ws.Range("F16").Value = 13
$ws.Range("F17").Value = 14
# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 882
$ws.Range("F3").Value = 4298
# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 882
$ws.Range("F3").Value = 35
$ws.Range("F6").Value = 1292
$ws.Range("F7").Value = 12995
$ws.Range("F8").Value = 287
$ws.Range("F9").Value = 741
$ws.Range("F10").Value = 4298
$ws.Range("F15").Value = 1887
$ws.Range("F16").Value = 41
$ws.Range("F18").Value = 4833
$ws.Range("F19").Value = 528
$ws.Range("F21").Value = 213
$ws.Range("F22").Value = 167
$ws.Range("F23").Value = 167
$ws.Range("F24").Value = 19
$ws.Range("F26").Value = 19
$ws.Range("F29").Value = 365
$ws.Range("F30").Value = 354
$ws.Range("F33").Value = 305
$ws.Range("F34").Value = 139
$ws.Range("F35").Value = 133
$ws.Range("F38").Value = 226
$ws.Range("F41").Value = 266
$ws.Range("F42").Value = 1315
$ws.Range("F43").Value = 13
$ws.Range("F44").Value = 348
$ws.Range("F47").Value = 14
